$wb = $excel.ActiveWorkbook

$wsJob = $wb.Worksheets.Item("Job to Run")
$wsAll = $wb.Worksheets.Item("All")

# Update the batch-template filename referenced on the "Job to Run" sheet
# (switch from the auto GNPS run template to the test-for-Python-workflow template)
$wsJob.Range("E2").Value = "MZmine3_batch_params_LCMSMS_HE_for_Commandline_2024_8_test_for_Python_workflow.xml"

# Update the selections on each sheet, then make "Job to Run" the active tab
# (selecting a range also activates its sheet, so do "All" first, "Job to Run" last)
$wsAll.Range("A2:E2").Select() | Out-Null
$wsJob.Range("B10").Select() | Out-Null
$wsJob.Activate() | Out-Null
